# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 23:52"

# 2. Update Estados Unidos row (row 4)
$ws.Range("B4").Value = 844992
$ws.Range("C4").Value = 26248
$ws.Range("D4").Value = 83910
$ws.Range("E4").Value = 713652
$ws.Range("G4").Value = 2112
$ws.Range("H4").Value = 47430

# 3. Update Brasil row (row 14)
$ws.Range("D14").Value = 25318
$ws.Range("E14").Value = 17533

# 4. Insert updated Congo data right after Paraguay (row 127), pushing
#    Islas Feroe and Gabon down one row each (their own numeric data is
#    unchanged, only their row position shifts). Martinica (row 131)
#    stays untouched.
$ws.Range("A128").Value = "Congo"
$ws.Range("B128").Value = 186
$ws.Range("C128").Value = 21
$ws.Range("D128").Value = 16
$ws.Range("E128").Value = 164
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 6

$ws.Range("A129").Value = "Islas Feroe"
$ws.Range("B129").Value = 185
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 178
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 0

$ws.Range("A130").Value = "Gabon"
$ws.Range("B130").Value = 166
$ws.Range("C130").Value = 10
$ws.Range("D130").Value = 24
$ws.Range("E130").Value = 141
$ws.Range("F130").Value = 2
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 1
